$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.865625
$ws.Cells.Item(2, 3).Value = 0.05036487956651937
$ws.Cells.Item(3, 2).Value = 0.91171875
$ws.Cells.Item(3, 3).Value = 0.02981060004427956
$ws.Cells.Item(4, 2).Value = 0.8609375
$ws.Cells.Item(4, 3).Value = 0.0488890245223813
$ws.Cells.Item(5, 2).Value = 0.50078125
$ws.Cells.Item(5, 3).Value = 0.04631171979774666
$ws.Cells.Item(6, 2).Value = 0.91953125
$ws.Cells.Item(6, 3).Value = 0.07430504851202911
$ws.Cells.Item(7, 2).Value = 0.86796875
$ws.Cells.Item(7, 3).Value = 0.02548360379734389
$ws.Cells.Item(8, 2).Value = 0.7953125
$ws.Cells.Item(8, 3).Value = 0.03702498628460246
$ws.Cells.Item(9, 2).Value = 0.8
$ws.Cells.Item(9, 3).Value = 0.05796011559684815
$ws.Cells.Item(10, 2).Value = 0.95859375
$ws.Cells.Item(10, 3).Value = 0.01858652695839113
$ws.Cells.Item(11, 2).Value = 0.853125
$ws.Cells.Item(11, 3).Value = 0.067775407482176
$ws.Cells.Item(12, 2).Value = 0.8636114137295191
$ws.Cells.Item(12, 3).Value = 0.0517951208072918
$ws.Cells.Item(13, 2).Value = 0.9110256386483482
$ws.Cells.Item(13, 3).Value = 0.03030858880563073
$ws.Cells.Item(14, 2).Value = 0.8582828815892848
$ws.Cells.Item(14, 3).Value = 0.05177958285980839
$ws.Cells.Item(15, 2).Value = 0.4850990583433807
$ws.Cells.Item(15, 3).Value = 0.03702907266559296
$ws.Cells.Item(16, 2).Value = 0.9185068986806142
$ws.Cells.Item(16, 3).Value = 0.07543230980102132
$ws.Cells.Item(17, 2).Value = 0.8656635937100488
$ws.Cells.Item(17, 3).Value = 0.02748353941379713
$ws.Cells.Item(18, 2).Value = 0.7807092927776933
$ws.Cells.Item(18, 3).Value = 0.04632330621588953
$ws.Cells.Item(19, 2).Value = 0.7752314982174235
$ws.Cells.Item(19, 3).Value = 0.08767573598631108
$ws.Cells.Item(20, 2).Value = 0.9583769564759038
$ws.Cells.Item(20, 3).Value = 0.01897274130644477
$ws.Cells.Item(21, 2).Value = 0.8505828951995115
$ws.Cells.Item(21, 3).Value = 0.06919544220630627
$ws.Cells.Item(22, 2).Value = 0.94879150390625
$ws.Cells.Item(22, 3).Value = 0.03869328532308846
$ws.Cells.Item(23, 2).Value = 0.9736083984375
$ws.Cells.Item(23, 3).Value = 0.02417946367052349
$ws.Cells.Item(24, 2).Value = 0.93121337890625
$ws.Cells.Item(24, 3).Value = 0.0260944530826347
$ws.Cells.Item(25, 2).Value = 0.50902099609375
$ws.Cells.Item(25, 3).Value = 0.0656675246379526
$ws.Cells.Item(26, 2).Value = 0.986767578125
$ws.Cells.Item(26, 3).Value = 0.02465875908253204
$ws.Cells.Item(27, 2).Value = 0.95079345703125
$ws.Cells.Item(27, 3).Value = 0.04838541984215908
$ws.Cells.Item(28, 2).Value = 0.951025390625
$ws.Cells.Item(28, 3).Value = 0.04990334705755102
$ws.Cells.Item(29, 2).Value = 0.9190185546875
$ws.Cells.Item(29, 3).Value = 0.035887642151469
$ws.Cells.Item(30, 2).Value = 0.99442138671875
$ws.Cells.Item(30, 3).Value = 0.002847563146840495
$ws.Cells.Item(31, 2).Value = 0.95325927734375
$ws.Cells.Item(31, 3).Value = 0.03112659885329115
